# The "Förändrad" (Changed) column (C) for every data row was bumped by one
# day, from 2023-10-03 (serial 45202) to 2023-10-04 (serial 45203).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$ws.Range("C2:C" + $lastRow).Value = 45203
